# Sync automático del tracker - update settled predictions (rows 220-248)
# and append a new pending prediction row (282) to the Predictions sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Predictions")

# Each tuple: Row, Status(L), Result(M), Resultado_Real(N), Profit(O), ROI(P), Enviado(Q)
$updates = @(
    @(220, "Completed", "Home Win", "Acierto", 1.2, 100, "2025-10-05 04:25:11"),
    @(221, "Completed", "Home Win", "Acierto", 1.43, 57, "2025-10-05 04:25:11"),
    @(222, "Completed", "Away Win", "Acierto", 0.34, 115, "2025-10-05 04:25:11"),
    @(223, "Completed", "Draw", "Fallo", -0.4, -100, "2025-10-05 04:25:11"),
    @(224, "Completed", "Draw", "Fallo", -2.1, -100, "2025-10-05 04:25:11"),
    @(225, "Completed", "Away Win", "Fallo", -1.7, -100, "2025-10-05 04:25:11"),
    @(226, "Completed", "Home Win", "Acierto", 1.61, 70, "2025-10-05 04:25:11"),
    @(227, "Completed", "Draw", "Fallo", -2.4, -100, "2025-10-05 04:25:11"),
    @(228, "Completed", "Home Win", "Acierto", 1.4, 50, "2025-10-05 04:25:11"),
    @(229, "Completed", "Away Win", "Fallo", -0.5, -100, "2025-10-05 04:25:11"),
    @(230, "Completed", "Home Win", "Acierto", 1.03, 115, "2025-10-05 04:25:11"),
    @(231, "Completed", "Draw", "Fallo", -2, -100, "2025-10-05 04:25:11"),
    @(232, "Completed", "Home Win", "Acierto", 1.43, 65, "2025-10-05 04:25:11"),
    @(233, "Completed", "Away Win", "Acierto", 1.26, 45, "2025-10-05 04:25:11"),
    @(234, "Completed", "Away Win", "Acierto", 1.18, 42, "2025-10-05 04:25:11"),
    @(235, "Completed", "Home Win", "Acierto", 1.48, 53, "2025-10-05 04:25:11"),
    @(236, "Completed", "Home Win", "Acierto", 1.45, 85, "2025-10-05 04:25:11"),
    @(237, "Completed", "Draw", "Fallo", -1.8, -100, "2025-10-05 04:25:11"),
    @(238, "Completed", "Draw", "Fallo", -2.7, -100, "2025-10-05 04:25:11"),
    @(239, "Completed", "Away Win", "Acierto", 1.12, 40, "2025-10-05 04:25:11"),
    @(240, "Completed", "Draw", "Fallo", -1.2, -100, "2025-10-05 15:19:34"),
    @(241, "Completed", "Draw", "Fallo", -0.8, -100, "2025-10-05 15:19:34"),
    @(242, "Completed", "Home Win", "Acierto", 1.62, 60, "2025-10-05 15:19:34"),
    @(243, "Completed", "Home Win", "Acierto", 1.26, 45, "2025-10-05 15:19:34"),
    @(244, "Completed", "Away Win", "Acierto", 1.23, 95, "2025-10-05 15:19:34"),
    @(245, "Completed", "Home Win", "Acierto", 1.58, 75, "2025-10-05 15:19:34"),
    @(246, "Completed", "Home Win", "Acierto", 1.18, 42, "2025-10-05 15:19:34"),
    @(247, "Completed", "Draw", "Fallo", -1.4, -100, "2025-10-05 15:19:34"),
    @(248, "Completed", "Home Win", "Acierto", 1.53, 90, "2025-10-05 15:19:34")
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 12).Value = $u[1]   # L - Status
    $ws.Cells.Item($r, 13).Value = $u[2]   # M - Result
    $ws.Cells.Item($r, 14).Value = $u[3]   # N - Resultado_Real
    $ws.Cells.Item($r, 15).Value = $u[4]   # O - Profit
    $ws.Cells.Item($r, 16).Value = $u[5]   # P - ROI
    $ws.Cells.Item($r, 17).Value = $u[6]   # Q - Enviado
}

# Append new row 282 with the latest pending prediction.
# Date-looking / percent-looking text must be forced to stay literal text
# (matching the source export's inline-string cells) instead of being
# auto-converted to a date serial / fraction by Excel's input parser:
# format the cell as Text first, assign the literal, then drop back to the
# workbook's default ("Normal") style so no stray number format lingers.
$newRow = 282

$cA = $ws.Cells.Item($newRow, 1)
$cA.NumberFormat = "@"
$cA.Value = "2025-10-06"
$cA.Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "Major League Soccer"
$ws.Cells.Item($newRow, 3).Value = "Los Angeles FC"
$ws.Cells.Item($newRow, 4).Value = "Atlanta United FC"
$ws.Cells.Item($newRow, 5).Value = "Home Win"

$cF = $ws.Cells.Item($newRow, 6)
$cF.NumberFormat = "@"
$cF.Value = "90.09%"
$cF.Style = "Normal"

$ws.Cells.Item($newRow, 7).Value = 1.4

$cH = $ws.Cells.Item($newRow, 8)
$cH.NumberFormat = "@"
$cH.Value = "24.87%"
$cH.Style = "Normal"

$ws.Cells.Item($newRow, 9).Value = 3.3
$ws.Cells.Item($newRow, 10).Value = 0.05
$ws.Cells.Item($newRow, 11).Value = 0.6533065723213444
$ws.Cells.Item($newRow, 12).Value = "Pending"

# M-Q stay blank (no result yet) - still materialise the cells themselves
# (format as text first) so the row shape matches the other data rows.
foreach ($col in 13..17) {
    $c = $ws.Cells.Item($newRow, $col)
    $c.NumberFormat = "@"
    $c.Value = ""
    $c.Style = "Normal"
}
